$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.509.87"
$ws.Range("E2").Value = "  +0.55%  "
$ws.Range("D3").Value = "1.915.67"
$ws.Range("E3").Value = "  +0.54%  "
$ws.Range("E4").Value = "  +0.61%  "
$ws.Range("D5").Value = "'325.83"
$ws.Range("E5").Value = "  +0.89%  "
$ws.Range("D6").Value = "'1.006"
$ws.Range("E6").Value = "  +0.54%  "
$ws.Range("D7").Value = "'0.4848"
$ws.Range("E7").Value = "  +2.54%  "
$ws.Range("D8").Value = "'0.4076"
$ws.Range("E8").Value = "  +1.04%  "
$ws.Range("D9").Value = "'0.08177"
$ws.Range("E9").Value = "  +2.06%  "
$ws.Range("D10").Value = "'1.016"
$ws.Range("E10").Value = "  +2.31%  "
$ws.Range("D11").Value = "'23.80"
$ws.Range("E11").Value = "  +5.18%  "
$ws.Range("D12").Value = "1.919.78"
$ws.Range("E12").Value = "  +3.80%  "
$ws.Range("D13").Value = "'6.050"
$ws.Range("E13").Value = "  +3.13%  "
$ws.Range("D14").Value = "'7.203"
$ws.Range("E14").Value = "  +2.01%  "
$ws.Range("D15").Value = "'91.18"
$ws.Range("E15").Value = "  +2.11%  "
$ws.Range("D16").Value = "'0.06777"
$ws.Range("E16").Value = "  +2.53%  "
$ws.Range("E17").Value = "  +0.58%  "
$ws.Range("D18").Value = "'0.00001041"
$ws.Range("E18").Value = "  +1.24%  "
$ws.Range("D19").Value = "'17.78"
$ws.Range("E19").Value = "  +1.28%  "
$ws.Range("D20").Value = "'1.005"
$ws.Range("E20").Value = "  +0.50%  "
$ws.Range("D21").Value = "29.531.42"
$ws.Range("E21").Value = "  +0.63%  "
$ws.Range("D22").Value = "'5.637"
$ws.Range("E22").Value = "  +2.24%  "
$ws.Range("E23").Value = "  +2.59%  "
$ws.Range("D24").Value = "'2.180"
$ws.Range("E24").Value = "  -0.92%  "
$ws.Range("D25").Value = "2.136.23"
$ws.Range("E25").Value = "  +1.95%  "
$ws.Range("D26").Value = "'156.51"
$ws.Range("E26").Value = "  +1.46%  "
$ws.Range("D27").Value = "'6.536"
$ws.Range("E27").Value = "  +8.25%  "
$ws.Range("E28").Value = "  +2.13%  "
$ws.Range("D29").Value = "'2.129"
$ws.Range("E29").Value = "  +1.69%  "
$ws.Range("D30").Value = "'120.74"
$ws.Range("E30").Value = "  +2.68%  "
$ws.Range("D31").Value = "'1.030"
$ws.Range("E31").Value = "  -3.34%  "
$ws.Range("D32").Value = "'0.09548"
$ws.Range("E32").Value = "  +0.85%  "
$ws.Range("D33").Value = "'5.521"
$ws.Range("E33").Value = "  +3.14%  "
$ws.Range("D34").Value = "'1.398"
$ws.Range("E34").Value = "  -1.04%  "
$ws.Range("D35").Value = "'3.559"
$ws.Range("E35").Value = "  +0.16%  "
$ws.Range("D36").Value = "'0.02281"
$ws.Range("E36").Value = "  +1.50%  "
$ws.Range("D37").Value = "'0.06135"
$ws.Range("E37").Value = "  +0.85%  "
$ws.Range("D38").Value = "'1.189"
$ws.Range("E39").Value = "  +7.99%  "
$ws.Range("D40").Value = "'0.5980"
$ws.Range("D41").Value = "'8.053"
$ws.Range("E41").Value = "  -0.10%  "
$ws.Range("D42").Value = "'0.1857"
$ws.Range("E42").Value = "  +1.59%  "
$ws.Range("D43").Value = "'2.410"
$ws.Range("E43").Value = "  -3.95%  "
$ws.Range("E44").Value = "  +0.64%  "
$ws.Range("B45").Value = "EnergySwap"
$ws.Range("C45").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D45").Value = "'12.52"
$ws.Range("E45").Value = "  +2.67%  "
$ws.Range("B46").Value = "Cronos"
$ws.Range("C46").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D46").Value = "'0.07632"
$ws.Range("E46").Value = "  -0.69%  "
$ws.Range("D47").Value = "'0.5592"
$ws.Range("D48").Value = "'1.959"
$ws.Range("E48").Value = "  +2.88%  "
$ws.Range("D49").Value = "'116.75"
$ws.Range("E49").Value = "  +2.84%  "
$ws.Range("D50").Value = "'72.79"
$ws.Range("E50").Value = "  +2.31%  "
$ws.Range("D51").Value = "'2.415"
$ws.Range("E51").Value = "  +3.15%  "
